# Rename the inline picture shapes in the headers/footers:
#   - the two Pearson Edexcel logo pictures (footers)   : image1.png -> image2.png
#   - the BTEC logo picture (header)                    : image2.jpg -> image1.jpg
#
# InlineShape has no writable .Name in the Word object model, so each
# picture is briefly promoted to a floating Shape (ConvertToShape),
# renamed there (Shape.Name -> wp:docPr/@name), and converted back to an
# inline picture (ConvertToInlineShape) in place.

$d = $word.ActiveDocument

function Rename-InlinePicture($inlineShape, $newName) {
    $shp = $inlineShape.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

foreach ($sec in $d.Sections) {

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $ishp = $shapes.Item($i)
                if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Rename-InlinePicture $ishp "image2.png"
                }
            }
        }
    }

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $ishp = $shapes.Item($i)
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlinePicture $ishp "image1.jpg"
                }
            }
        }
    }
}
